$d = $word.ActiveDocument

# --- Change 1: body text -----------------------------------------------
# After the paragraph "Conceitos de Qualidade de Software. " insert two
# blank paragraphs followed by a paragraph containing the new text
# "Teste 1 – alterando documento..", all sharing the paragraph
# formatting (justify both, sz/szCs 28) of the paragraph that follows.

$anchorIndex = 0
$i = 0
foreach ($p in $d.Paragraphs) {
    $i = $i + 1
    if ($p.Range.Text.Contains("Conceitos de Qualidade de Software")) {
        $anchorIndex = $i
    }
}

if ($anchorIndex -gt 0) {
    $nextPara = $d.Paragraphs.Item($anchorIndex + 1)
    $insertRange = $nextPara.Range

    $insertRange.InsertParagraphBefore()
    $insertRange.InsertParagraphBefore()
    $insertRange.InsertParagraphBefore()

    $newTextPara = $d.Paragraphs.Item($anchorIndex + 3)
    $newTextPara.Range.InsertAfter("Teste 1 – alterando documento..")
}

# --- Change 2: footer ----------------------------------------------------
# Collapse the three runs "Prof. " / "Keity" / " Yamamoto" (with the
# spell-check proofErr markers around "Keity") into a single run.

$sec = $d.Sections.Item(1)
$ftr = $sec.Footers.Item(1)
$ftrRange = $ftr.Range
$ftrRange.Find.Execute("Prof. Keity Yamamoto", $true, $false, $false, $false, $false, $true, 1, $false, "Prof. Keity Yamamoto", 2)

Write-Output "done"
